$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

# Row 2
Set-TextValue 2 4 "28.217.52"
$ws.Cells.Item(2, 5).Value = "  -2.64%  "

# Row 3
Set-TextValue 3 4 "1.866.18"
$ws.Cells.Item(3, 5).Value = "  -2.28%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.19%  "

# Row 5
Set-TextValue 5 4 "318.85"
$ws.Cells.Item(5, 5).Value = "  -1.88%  "

# Row 6
$ws.Cells.Item(6, 5).Value = "  +0.22%  "

# Row 7
Set-TextValue 7 4 "0.4391"
$ws.Cells.Item(7, 5).Value = "  -4.35%  "

# Row 8
Set-TextValue 8 4 "0.3688"
$ws.Cells.Item(8, 5).Value = "  -3.54%  "

# Row 9
Set-TextValue 9 4 "0.07487"
$ws.Cells.Item(9, 5).Value = "  -3.00%  "

# Row 10
Set-TextValue 10 4 "0.9347"
$ws.Cells.Item(10, 5).Value = "  -4.69%  "

# Row 11
Set-TextValue 11 4 "21.28"
$ws.Cells.Item(11, 5).Value = "  -3.63%  "

# Row 12
Set-TextValue 12 4 "1.978.18"
$ws.Cells.Item(12, 5).Value = "  +3.72%  "

# Row 13
Set-TextValue 13 4 "6.699"
$ws.Cells.Item(13, 5).Value = "  -3.40%  "

# Row 14
Set-TextValue 14 4 "5.467"
$ws.Cells.Item(14, 5).Value = "  -3.70%  "

# Row 15
Set-TextValue 15 4 "0.06906"
$ws.Cells.Item(15, 5).Value = "  -1.63%  "

# Row 16
$ws.Cells.Item(16, 5).Value = "  +0.23%  "

# Row 17
Set-TextValue 17 4 "81.78"
$ws.Cells.Item(17, 5).Value = "  -2.56%  "

# Row 18
Set-TextValue 18 4 "0.000009019"
$ws.Cells.Item(18, 5).Value = "  -4.77%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  +0.15%  "

# Row 20
Set-TextValue 20 4 "15.88"
$ws.Cells.Item(20, 5).Value = "  -5.08%  "

# Row 21
Set-TextValue 21 4 "28.204.67"
$ws.Cells.Item(21, 5).Value = "  -2.54%  "

# Row 22
Set-TextValue 22 4 "5.114"
$ws.Cells.Item(22, 5).Value = "  -3.80%  "

# Row 23
Set-TextValue 23 4 "10.78"
$ws.Cells.Item(23, 5).Value = "  -0.96%  "

# Row 24
Set-TextValue 24 4 "2.129.04"
$ws.Cells.Item(24, 5).Value = "  -1.09%  "

# Row 25
Set-TextValue 25 4 "2.021"
$ws.Cells.Item(25, 5).Value = "  -3.40%  "

# Row 26
Set-TextValue 26 4 "155.12"
$ws.Cells.Item(26, 5).Value = "  -2.06%  "

# Row 27
Set-TextValue 27 4 "18.33"
$ws.Cells.Item(27, 5).Value = "  -3.72%  "

# Row 28
Set-TextValue 28 4 "5.316"
$ws.Cells.Item(28, 5).Value = "  -6.18%  "

# Row 29
Set-TextValue 29 4 "113.28"
$ws.Cells.Item(29, 5).Value = "  -3.78%  "

# Row 30
Set-TextValue 30 4 "1.722"
$ws.Cells.Item(30, 5).Value = "  -6.95%  "

# Row 31
Set-TextValue 31 4 "0.09008"
$ws.Cells.Item(31, 5).Value = "  -2.84%  "

# Row 32
$ws.Cells.Item(32, 2).Value = "Filecoin"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue 32 4 "4.836"
$ws.Cells.Item(32, 5).Value = "  -4.64%  "

# Row 33
$ws.Cells.Item(33, 2).Value = "ImmutableX"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue 33 4 "0.7913"
$ws.Cells.Item(33, 5).Value = "  -8.56%  "

# Row 34
Set-TextValue 34 4 "1.173"
$ws.Cells.Item(34, 5).Value = "  -6.21%  "

# Row 35
Set-TextValue 35 4 "2.940"
$ws.Cells.Item(35, 5).Value = "  -2.80%  "

# Row 37
Set-TextValue 37 4 "1.125"
$ws.Cells.Item(37, 5).Value = "  -2.58%  "

# Row 38
Set-TextValue 38 4 "0.05433"
$ws.Cells.Item(38, 5).Value = "  -5.51%  "

# Row 39
Set-TextValue 39 4 "0.01967"
$ws.Cells.Item(39, 5).Value = "  -3.61%  "

# Row 40
Set-TextValue 40 4 "2.952"
$ws.Cells.Item(40, 5).Value = "  +2.73%  "

# Row 41
Set-TextValue 41 4 "0.5251"
$ws.Cells.Item(41, 5).Value = "  -4.85%  "

# Row 42
Set-TextValue 42 4 "6.993"
$ws.Cells.Item(42, 5).Value = "  -5.91%  "

# Row 43
Set-TextValue 43 4 "0.1677"
$ws.Cells.Item(43, 5).Value = "  -4.58%  "

# Row 44
Set-TextValue 44 4 "8.704"
$ws.Cells.Item(44, 5).Value = "  -6.69%  "

# Row 45
Set-TextValue 45 4 "0.06738"
$ws.Cells.Item(45, 5).Value = "  -1.44%  "

# Row 46
Set-TextValue 46 4 "0.4863"
$ws.Cells.Item(46, 5).Value = "  -6.06%  "

# Row 47
Set-TextValue 47 4 "10.57"
$ws.Cells.Item(47, 5).Value = "  -5.90%  "

# Row 48
Set-TextValue 48 4 "106.81"
$ws.Cells.Item(48, 5).Value = "  -3.77%  "

# Row 49
Set-TextValue 49 4 "1.922"
$ws.Cells.Item(49, 5).Value = "  -6.80%  "

# Row 50
$ws.Cells.Item(50, 5).Value = "  +0.18%  "

# Row 51
Set-TextValue 51 4 "1.669"
$ws.Cells.Item(51, 5).Value = "  -6.42%  "
